$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 160341
$ws.Range("C4").Value = 151376
$ws.Range("C5").Value = 8966
$ws.Range("C8").Value = 64.44
